$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Bus index" sheet: wind-trace name for bus 1 changed (NEN -> HUN)
# ---------------------------------------------------------------------------
$wsBusIndex = $wb.Worksheets.Item("Bus index")
$wsBusIndex.Range("D2").Value = "Bubble HUN Wind 2014-2045_0910refyr.csv"

# ---------------------------------------------------------------------------
# "Bus connections" sheet: several reactance formulas / MW limits updated
# ---------------------------------------------------------------------------
$wsBusConn = $wb.Worksheets.Item("Bus connections")

$wsBusConn.Range("C2").Formula = "=(PI()/12)/(E2/100)*D2"
$wsBusConn.Range("C3").Formula = "=(PI()/12)/(E3/100)*D3"

$wsBusConn.Range("C4").Formula = "=(PI()/6)/(E4/100)*D4"
$wsBusConn.Range("E4").Value = 2700

$wsBusConn.Range("C5").Formula = "=(PI()/3)/(E5/100)*D5"
$wsBusConn.Range("E5").Value = 1800

$wsBusConn.Range("C6").Formula = "=(PI()/6)/(E6/100)*D6"
$wsBusConn.Range("E6").Value = 4800

$wsBusConn.Range("C7").Formula = "=(PI()/6)/(E7/100)*D7"

$wsBusConn.Range("C8").Formula = "=(PI()/6)/(E8/100)*D8"
$wsBusConn.Range("E8").Value = 4600

# ---------------------------------------------------------------------------
# View-state bookkeeping: which sheet/cell was selected when the file was
# last saved.
# ---------------------------------------------------------------------------
$wsBusConn.Range("A2:D8").Select()

$wsBusIndex.Activate()
$wsBusIndex.Range("D8").Select()
